# Append new scrape run: 2025-12-03 12:39 JST
# Prepend a new job listing, update all existing timestamps, and insert one
# more new listing right before the last (existing) row, matching the
# source site's ordering for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ts = "2025-12-03 12:39:22"

# --- Shift existing data down, opening two new rows -----------------------
# Row 2: brand-new top listing for this run.
$ws.Rows.Item(2).Insert()
# Row 9 (old row 8 is now at row 9 after the first insert): another new
# listing the site added between the previous run's rows 7 and 8.
$ws.Rows.Item(9).Insert()

# --- Refresh the timestamp on every (now-shifted) existing data row -------
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = $ts
}
$ws.Cells.Item(10, 1).Value = $ts

# --- New row 2 ---------------------------------------------------------
$ws.Cells.Item(2, 1).Value = $ts
$ws.Cells.Item(2, 2).Value = "【急募】LLMを活用した次世代AIエージェント開発プロジェクト"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5446583"
$ws.Cells.Item(2, 7).Value = 375
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発"

# --- New row 9 ---------------------------------------------------------
$ws.Cells.Item(9, 1).Value = $ts
$ws.Cells.Item(9, 2).Value = "マッチングアプリ(Web)のPMできる方を募集します!"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5446673"
$ws.Cells.Item(9, 7).Value = 45
$ws.Cells.Item(9, 8).Value = "◇アプリ"

# --- Rebuild the hyperlinks on column F, rows 2..10 ------------------------
# Row-insert doesn't renumber the existing hyperlink relationships, so drop
# them all and re-add fresh ones against the final row layout.
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5446583",
    "https://www.lancers.jp/work/detail/5446091",
    "https://www.lancers.jp/work/detail/5446106",
    "https://www.lancers.jp/work/detail/5446360",
    "https://www.lancers.jp/work/detail/5446289",
    "https://www.lancers.jp/work/detail/5446279",
    "https://www.lancers.jp/work/detail/5446285",
    "https://www.lancers.jp/work/detail/5446673",
    "https://www.lancers.jp/work/detail/5446233"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i])
}

# Re-apply the built-in Hyperlink cell style so every linked cell reuses the
# same style record the sheet already had (Hyperlinks.Add otherwise mints a
# redundant duplicate style entry).
$ws.Range("F2:F10").Style = "Hyperlink"
